$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.07460170201535693
$ws.Range("C2").Value = 0.4066704939480233
$ws.Range("D2").Value = 0.2229059984459623
$ws.Range("E2").Value = 0.472129217954113
$ws.Range("F2").Value = 0.4825607857740096
$ws.Range("G2").Value = 15

$ws.Range("B3").Value = 0.314504379002312
$ws.Range("C3").Value = 0.4790566883332484
$ws.Range("D3").Value = 0.3549208902847111
$ws.Range("E3").Value = 0.5957523732933937
$ws.Range("F3").Value = 0.5250721856764178
$ws.Range("G3").Value = 14

$ws.Range("B4").Value = 0.4632337633744819
$ws.Range("C4").Value = 0.5692069213481521
$ws.Range("D4").Value = 0.4862027082318384
$ws.Range("E4").Value = 0.6972823733838669
$ws.Range("F4").Value = 0.5424499556889087
$ws.Range("G4").Value = 13

$ws.Range("B5").Value = 0.6092340853067807
$ws.Range("C5").Value = 0.6750741196220317
$ws.Range("D5").Value = 0.6306383360763419
$ws.Range("E5").Value = 0.7941274054434476
$ws.Range("F5").Value = 0.5320343448005646
$ws.Range("G5").Value = 12

$ws.Range("B6").Value = 0.638688239432181
$ws.Range("C6").Value = 0.7242386063103117
$ws.Range("D6").Value = 0.7271270905885431
$ws.Range("E6").Value = 0.8527174740724756
$ws.Range("F6").Value = 0.5925579007485432
$ws.Range("G6").Value = 11

$ws.Range("B7").Value = 0.5167424503166
$ws.Range("C7").Value = 0.6241159593079864
$ws.Range("D7").Value = 0.538927811949409
$ws.Range("E7").Value = 0.7341170287831559
$ws.Range("F7").Value = 0.5496514572286344
$ws.Range("G7").Value = 10

$ws.Range("B8").Value = 0.370868480679851
$ws.Range("C8").Value = 0.5185686719842275
$ws.Range("D8").Value = 0.3627064250084273
$ws.Range("E8").Value = 0.6022511311806955
$ws.Range("F8").Value = 0.5032974959479504
$ws.Range("G8").Value = 9

$ws.Range("B9").Value = 0.2966541511547854
$ws.Range("C9").Value = 0.4816046978711712
$ws.Range("D9").Value = 0.2612550586842152
$ws.Range("E9").Value = 0.5111311560492231
$ws.Range("F9").Value = 0.4559623317163587
$ws.Range("G9").Value = 6

$ws.Range("B10").Value = 0.4961513153058069
$ws.Range("C10").Value = 0.4961513153058069
$ws.Range("D10").Value = 0.3214734301341839
$ws.Range("E10").Value = 0.5669862697933556
$ws.Range("F10").Value = 0.3360966433657923
$ws.Range("G10").Value = 3

